$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing existing rows 16-35 down to 17-36.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with this week's reading
# (same as the last entry in the series, one reporting cycle later).
$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C16").Value = "Los Lagos"
$ws.Range("D16").Value2 = 44915
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100101
$ws.Range("H16").Value = "Berries"
$ws.Range("I16").Value = 100101001
$ws.Range("J16").Value = "Arándano (blue)"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 600
$ws.Range("N16").Value = 4000
$ws.Range("O16").Value = 4200
$ws.Range("P16").Value = 4100
$ws.Range("Q16").Value = "$/bandeja 2 kilos"
$ws.Range("R16").Value = "Provincia de Curicó"
$ws.Range("S16").Value = 2050
$ws.Range("T16").Value = 2
